$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1649.9615
$ws.Range("J17").Value = 1649.9615
$ws.Range("L17").Value = 4949.8845
$ws.Range("N17").Value = -5285.8845
$ws.Range("H55").Value = 529
$ws.Range("J55").Value = 773.5
$ws.Range("L55").Value = 773.5
$ws.Range("N55").Value = -1201.5
$ws.Range("H107").Value = 1625
$ws.Range("I107").Value = 1333.3334
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 1333.3334
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 586.6666
$ws.Range("N107").Value = -5640
$ws.Range("H125").Value = 2867
$ws.Range("I125").Value = 2263.7778
$ws.Range("J125").Value = 4224.25
$ws.Range("K125").Value = 20374.0002
$ws.Range("L125").Value = 38018.25
$ws.Range("M125").Value = -17914.0002
$ws.Range("N125").Value = -42938.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 15402.571
$ws.Range("I16").Value = 33624
$ws.Range("J16").Value = 1736.5
$ws.Range("K16").Value = 33624
$ws.Range("L16").Value = 1736.5
$ws.Range("M16").Value = -33337
$ws.Range("N16").Value = -2310.5
$ws.Range("H74").Value = 1116.0714
$ws.Range("I74").Value = 1052.7273
$ws.Range("J74").Value = 1348.3334
$ws.Range("K74").Value = 1052.7273
$ws.Range("L74").Value = 1348.3334
$ws.Range("M74").Value = -178.7273
$ws.Range("N74").Value = -3096.3334
$ws.Range("H77").Value = 1116.0714
$ws.Range("I77").Value = 1052.7273
$ws.Range("J77").Value = 1348.3334
$ws.Range("K77").Value = 5263.636500000001
$ws.Range("L77").Value = 6741.666999999999
$ws.Range("M77").Value = -895.6365000000005
$ws.Range("N77").Value = -15477.667
$ws.Range("H88").Value = 2175.4614
$ws.Range("J88").Value = 2372
$ws.Range("L88").Value = 2372
$ws.Range("N88").Value = -3184
$ws.Range("H91").Value = 2175.4614
$ws.Range("J91").Value = 2372
$ws.Range("L91").Value = 2372
$ws.Range("N91").Value = -5180

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 78000
$ws.Range("J74").Value = 78000
$ws.Range("L74").Value = 78000
$ws.Range("N74").Value = -79872
$ws.Range("H77").Value = 78000
$ws.Range("J77").Value = 78000
$ws.Range("L77").Value = 234000
$ws.Range("N77").Value = -243360
$ws.Range("H99").Value = 2382.1667
$ws.Range("I99").Value = 2323.5
$ws.Range("K99").Value = 2323.5
$ws.Range("M99").Value = -825.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 16981
$ws.Range("J17").Value = 16981
$ws.Range("L17").Value = 16981
$ws.Range("N17").Value = -17329
$ws.Range("H132").Value = 2579.4
$ws.Range("I132").Value = 2421.5557
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7264.6671
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4734.6671
$ws.Range("N132").Value = -17060
$ws.Range("H134").Value = 1843.2
$ws.Range("I134").Value = 1791.5
$ws.Range("J134").Value = 2050
$ws.Range("K134").Value = 5374.5
$ws.Range("L134").Value = 6150
$ws.Range("M134").Value = -2839.5
$ws.Range("N134").Value = -11220

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 237.66667
$ws.Range("J2").Value = 266.2
$ws.Range("L2").Value = 1597.2
$ws.Range("N2").Value = -1823.2
$ws.Range("H46").Value = 2333.5833
$ws.Range("I46").Value = 2237.5
$ws.Range("J46").Value = 2352.8
$ws.Range("K46").Value = 6712.5
$ws.Range("L46").Value = 7058.400000000001
$ws.Range("M46").Value = -6621.5
$ws.Range("N46").Value = -7240.400000000001
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H122").Value = 1377.8
$ws.Range("I122").Value = 697
$ws.Range("J122").Value = 1831.6666
$ws.Range("K122").Value = 6273
$ws.Range("L122").Value = 16484.9994
$ws.Range("N122").Value = -21384.9994
$ws.Range("M122").Value = -3823
$ws.Range("H131").Value = 765
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 3998.4285
$ws.Range("I132").Value = 2664.8333
$ws.Range("K132").Value = 23983.4997
$ws.Range("M132").Value = -21453.4997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 29790
$ws.Range("J62").Value = 29790
$ws.Range("L62").Value = 29790
$ws.Range("N62").Value = -31162
$ws.Range("H65").Value = 29790
$ws.Range("J65").Value = 29790
$ws.Range("L65").Value = 89370
$ws.Range("N65").Value = -96234
$ws.Range("H102").Value = 2944.818
$ws.Range("I102").Value = 1710.6666
$ws.Range("K102").Value = 1710.6666
$ws.Range("M102").Value = -88.66660000000002
$ws.Range("H135").Value = 49000
$ws.Range("J135").Value = 49000
$ws.Range("L135").Value = 49000
$ws.Range("N135").Value = -59140
$ws.Range("H140").Value = 69420
$ws.Range("J140").Value = 69420
$ws.Range("L140").Value = 69420
$ws.Range("N140").Value = -79780

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1435.75
$ws.Range("I22").Value = 1099.4
$ws.Range("K22").Value = 1099.4
$ws.Range("M22").Value = -804.4000000000001
$ws.Range("H27").Value = 1435.75
$ws.Range("I27").Value = 1099.4
$ws.Range("K27").Value = 1099.4
$ws.Range("M27").Value = -992.4000000000001
$ws.Range("H46").Value = 3789.4211
$ws.Range("I46").Value = 2000
$ws.Range("K46").Value = 2000
$ws.Range("M46").Value = -1812
$ws.Range("H55").Value = 195.42857
$ws.Range("I55").Value = 190.66667
$ws.Range("K55").Value = 190.66667
$ws.Range("M55").Value = -17.66667000000001
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 16315.223
$ws.Range("I132").Value = 10571
$ws.Range("J132").Value = 19970.637
$ws.Range("K132").Value = 31713
$ws.Range("L132").Value = 59911.91099999999
$ws.Range("M132").Value = -29183
$ws.Range("N132").Value = -64971.91099999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1197.1
$ws.Range("I100").Value = 1263.4445
$ws.Range("K100").Value = 2526.889
$ws.Range("M100").Value = -1985.889
$ws.Range("H113").Value = 562.63635
$ws.Range("I113").Value = 609.2222
$ws.Range("J113").Value = 353
$ws.Range("K113").Value = 1827.6666
$ws.Range("L113").Value = 1059
$ws.Range("M113").Value = 342.3334
$ws.Range("N113").Value = -5399
